# The deck's single design theme (ppt/theme/theme1.xml, used by the slide
# master / all slides) currently carries the "Integral" colour scheme.
# The commit swaps it for the stock "Office Theme" colour scheme (the
# theme that used to only back the notes master). Font scheme and format
# scheme (fills/lines/effects) are identical between the two themes in
# this deck, so only the 12 theme colours actually need to change.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

$colorScheme.Colors(1).RGB  = 0x000000   # Dark 1    -> 000000
$colorScheme.Colors(2).RGB  = 0xFFFFFF   # Light 1   -> FFFFFF
$colorScheme.Colors(3).RGB  = 0x6A5444   # Dark 2    -> 44546A
$colorScheme.Colors(4).RGB  = 0xE6E6E7   # Light 2   -> E7E6E6
$colorScheme.Colors(5).RGB  = 0xD59B5B   # Accent 1  -> 5B9BD5
$colorScheme.Colors(6).RGB  = 0x317DED   # Accent 2  -> ED7D31
$colorScheme.Colors(7).RGB  = 0xA5A5A5   # Accent 3  -> A5A5A5
$colorScheme.Colors(8).RGB  = 0x00C0FF   # Accent 4  -> FFC000
$colorScheme.Colors(9).RGB  = 0xC47244   # Accent 5  -> 4472C4
$colorScheme.Colors(10).RGB = 0x47AD70   # Accent 6  -> 70AD47
$colorScheme.Colors(11).RGB = 0xC16305   # Hyperlink -> 0563C1
$colorScheme.Colors(12).RGB = 0x724F95   # Followed Hyperlink -> 954F72
